$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sheet (tab) name to reflect the new "through" date.
$ws.Name = "Through 2021-09-14"

# Update the header label in A10 ("September (through 09-13)" -> "...09-14").
$ws.Cells.Item(10, 1).Value = "September (through 09-14)"

# March 2021 (column H, row 4) was revised down by one.
$ws.Cells.Item(4, 8).Value = 81

# September row (row 10) gets the new daily counts (F10 is unchanged).
$ws.Cells.Item(10, 2).Value = 15
$ws.Cells.Item(10, 3).Value = 24
$ws.Cells.Item(10, 4).Value = 35
$ws.Cells.Item(10, 5).Value = 27
$ws.Cells.Item(10, 7).Value = 53
$ws.Cells.Item(10, 8).Value = 70

# Total row (row 11) reflects the updated year-to-date sums (F11 is unchanged).
$ws.Cells.Item(11, 2).Value = 209
$ws.Cells.Item(11, 3).Value = 405
$ws.Cells.Item(11, 4).Value = 586
$ws.Cells.Item(11, 5).Value = 517
$ws.Cells.Item(11, 7).Value = 837
$ws.Cells.Item(11, 8).Value = 1140
